# Slide and Catch pt.2 updates
#
# The "Sprite assets" / "Non-Sprite assets" bullet lists used
# "coming in part 2" as a placeholder for items that hadn't been
# implemented yet. Those pieces are now done, so replace the
# placeholders with their real descriptions (or, where nothing more
# needs to be said, just drop the placeholder text).

$d = $word.ActiveDocument

function Replace-Once($oldText, $newText) {
    # Search the whole document, replacing only the first remaining
    # match. "coming in part 2" only ever appears as the tail of the
    # (non-bold) description run, so this never touches the bold
    # label name that precedes it, and only the matched characters are
    # swapped - the surrounding runs/formatting are left alone.
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute(
        $oldText,   # FindText
        $false,     # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $newText,   # ReplaceWith
        1           # Replace (wdReplaceOne)
    )
}

# lblLives - now implemented as a simpleGE Label showing the number of
# lives remaining.
Replace-Once "coming in part 2" "a simpleGE Label displaying the number of lives"

# lblTime - now implemented as a simpleGE Label showing the elapsed
# time.
Replace-Once "coming in part 2" "a simpleGE Label displaying time"

# lblScore - now implemented as a simpleGE Label showing the score
# expressed as time in seconds.
Replace-Once "coming in part 2" "a simpleGE Label displaying score as time in secs "

# timer, score, lives, sndAxe - these were only ever placeholders with
# no further description; just drop the "coming in part 2" text and
# leave the rest of the line (and its formatting/runs) untouched.
Replace-Once "coming in part 2" ""
Replace-Once "coming in part 2" ""
Replace-Once "coming in part 2" ""
Replace-Once "coming in part 2" ""
